$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 does not exist yet (gap between row4 and row6). Row 6 already carries the
# formatting (styles) that row 5 needs, so copy that formatting down into row 5
# before we fill row 6 with its own data. (Only A:J - row 5 has no K cell.)
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 5: new data row (Crumpet GEF / Crumpet exporter)
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6: fill in existing (previously empty) row with data (Scone GEF / Scone exporter)
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Update the selection/view state to match the final file
$ws.Range("D7").Select()
